$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Cells.Item(3, 2).Value = "6.0.0"

# Date changes
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# Delete the duplicate "Contact" row (row 10), shifting rows 11-21 up by one
$meta.Rows.Item(10).Delete()

# The row that is now row 10 (was "Contact"/"No display for ContactDetail") becomes Jurisdiction
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# --- "Elements" sheet (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 (the root Extension element): Short/Definition columns (K/L) updated
$elements.Cells.Item(2, 11).Value = "Drug Fully Insured Indicator"
$elements.Cells.Item(2, 12).Value = "Indicator of the fully insured drug coverage for the member or employee"
